$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (cell, newValue) updates for column G
# (Elapsed Duration(Hrs)) as described by the diff.
$updates = @{
    "R1" = @(
        @{ Cell = "G2"; Value = "3926:16:37" },
        @{ Cell = "G3"; Value = "65:49:15" }
    )
    "R2" = @(
        @{ Cell = "G2"; Value = "12107:40:16" },
        @{ Cell = "G3"; Value = "3237:23:45" },
        @{ Cell = "G4"; Value = "475:35:19" }
    )
    "R4" = @(
        @{ Cell = "G2"; Value = "2953:30:05" },
        @{ Cell = "G3"; Value = "180:42:20" }
    )
    "R5" = @(
        @{ Cell = "G2"; Value = "427:29:04" }
    )
    "R6" = @(
        @{ Cell = "G2"; Value = "68:01:22" }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($update in $updates[$sheetName]) {
        $ws.Range($update.Cell).Value = $update.Value
    }
}
